$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# The target slide already has shapes with ids 1,2,4,8 in use. PowerPoint's
# shape-id allocator hands out the next free integer (skipping ids already
# in use) and does not reuse ids once handed out. Adding + deleting one
# throw-away textbox first consumes id 3 so that the two real textboxes we
# add next land on ids 5 and 6, matching the target deck exactly.
$dummy = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$dummy.Delete()

# --- TextBox 4: "draft-ietf-opsawg-oam-characterization" -------------------
$tb1 = $s.Shapes.AddTextbox(1, 131.12574803149607, 268.51929133858266, 209.23811023622048, 19.38748031496063)
$tb1.Name = "TextBox 4"
$tb1.Fill.Visible = $false
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "draft-ietf-opsawg-oam-characterization"

# Five separate runs (word boundaries), each carrying identical character
# formatting, mirroring the way PowerPoint split them while the author
# was typing / spell-checking.
$segments1 = @(
    @{start=1;  len=6}   # "draft-"
    @{start=7;  len=4}   # "ietf"
    @{start=11; len=8}   # "-opsawg-"
    @{start=19; len=3}   # "oam"
    @{start=22; len=18}  # "-characterization"
)
foreach ($seg in $segments1) {
    $run = $tr1.Characters($seg.start, $seg.len)
    $run.LanguageID = "en-GB"
    $run.Font.NameFarEast = "Inter"
    $run.Font.Name = "Inter"
    $run.Font.Size = 10
    $run.Font.Color.ObjectThemeColor = 8
}

# --- TextBox 5: "T. Graf" ---------------------------------------------------
$tb2 = $s.Shapes.AddTextbox(1, 427.6806299212598, 268.5192125984252, 73.9651968503937, 19.38748031496063)
$tb2.Name = "TextBox 5"
$tb2.Fill.Visible = $false
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "T. Graf"
$tr2.LanguageID = "en-GB"
$tr2.Font.NameFarEast = "Inter"
$tr2.Font.Name = "Inter"
$tr2.Font.Size = 10
$tr2.Font.Color.ObjectThemeColor = 8
